$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("quiz")

# Update "Marking" row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (B12): 69 -> 115
$ws.Range("B12").Value = 115

# Update the correct/total marks label (E12): "68/84" -> "115/140"
$ws.Range("E12").Value = "115/140"
